$d = $word.ActiveDocument

# --- "Programa" section: split the single run's text into three text runs
# separated by manual line breaks (<w:br/>), right after
# "...condicional/marginal." and right after "...e independentes.".
# The anchor text below only occurs in the "Programa" paragraph (the
# "Programa resumido" paragraph uses lower-case "variância" instead of
# "Variância"), so this Find/Replace is unambiguous.

$found1 = $d.Content.Find.Execute(
    "Variância condicional/marginal.Estatística",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Variância condicional/marginal.^lEstatística", 2)
if (-not $found1) { throw "anchor 1 (Programa) not found" }

$found2 = $d.Content.Find.Execute(
    "amostras pareadas e independentes.Técnicas",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "amostras pareadas e independentes.^lTécnicas", 2)
if (-not $found2) { throw "anchor 2 (Programa) not found" }

# --- "Bibliografia" section: split the single run's text into four text
# runs separated by a blank line (two manual line breaks), one per
# reference entry.

$found3 = $d.Content.Find.Execute(
    "LTC, 2009. D. C. Montgomery",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "LTC, 2009. ^l^lD. C. Montgomery", 2)
if (-not $found3) { throw "anchor 3 (Bibliografia) not found" }

$found4 = $d.Content.Find.Execute(
    "John Wiley, 2006.W. J. Conover",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "John Wiley, 2006.^l^lW. J. Conover", 2)
if (-not $found4) { throw "anchor 4 (Bibliografia) not found" }

$found5 = $d.Content.Find.Execute(
    "John Wiley d Sons, 1999.R. A. Johnson",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "John Wiley d Sons, 1999.^l^lR. A. Johnson", 2)
if (-not $found5) { throw "anchor 5 (Bibliografia) not found" }

Write-Output "OK: $found1 $found2 $found3 $found4 $found5"
